# Update automàtic: dades i banners [2026-02-17 16:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Known-good unaffected cell in the same style class (s="3"),
# used as a format donor so risky coercible text (e.g. "75%")
# keeps its original style after being forced to literal text.
$formatDonor = $ws.Range("H2")

$ws.Range("E2").Value = "2026-02-17 16:18:49"
$ws.Range("K2").Value = "6.6 MJ/m2"
$ws.Range("E3").Value = "2026-02-17 16:18:51"
$ws.Range("K3").Value = "6.4 MJ/m2"
$ws.Range("O3").Value = "-4.7 °C"
$ws.Range("E4").Value = "2026-02-17 16:18:54"
$ws.Range("K4").Value = "7.2 MJ/m2"
$ws.Range("O4").Value = "9.4 °C"
$ws.Range("E5").Value = "2026-02-17 16:18:56"
$ws.Range("K5").Value = "5.6 MJ/m2"
$ws.Range("O5").Value = "-4.2 °C"
$ws.Range("E6").Value = "2026-02-17 16:18:59"
$ws.Range("H6").Value = "'75%"
$formatDonor.Copy() | Out-Null
$ws.Range("H6").PasteSpecial(-4122) | Out-Null
$ws.Range("K6").Value = "9.6 MJ/m2"
$ws.Range("O6").Value = "10.4 °C"
$ws.Range("E7").Value = "2026-02-17 16:19:02"
$ws.Range("K7").Value = "12.8 MJ/m2"
$ws.Range("E8").Value = "2026-02-17 16:19:04"
$ws.Range("J8").Value = "1017.8 hPa"
$ws.Range("K8").Value = "13.7 MJ/m2"
$ws.Range("E9").Value = "2026-02-17 16:19:07"
$ws.Range("K9").Value = "8.4 MJ/m2"
$ws.Range("O9").Value = "12.8 °C"
$ws.Range("E10").Value = "2026-02-17 16:19:09"
$ws.Range("K10").Value = "10.2 MJ/m2"
$ws.Range("O10").Value = "10.7 °C"
$ws.Range("E11").Value = "2026-02-17 16:19:12"
$ws.Range("O11").Value = "7.5 °C"
$ws.Range("E12").Value = "2026-02-17 16:19:15"
$ws.Range("E13").Value = "2026-02-17 16:19:17"
$ws.Range("J13").Value = "1017.6 hPa"
$ws.Range("K13").Value = "9.0 MJ/m2"
$ws.Range("O13").Value = "6.8 °C"
$ws.Range("E14").Value = "2026-02-17 16:19:20"
$ws.Range("K14").Value = "13.2 MJ/m2"
$ws.Range("O14").Value = "13.7 °C"
$ws.Range("E15").Value = "2026-02-17 16:19:22"
$ws.Range("E16").Value = "2026-02-17 16:19:25"
$ws.Range("K16").Value = "9.7 MJ/m2"
$ws.Range("O16").Value = "-3.9 °C"
$ws.Range("E17").Value = "2026-02-17 16:19:27"
$ws.Range("E18").Value = "2026-02-17 16:19:30"
$ws.Range("H18").Value = "'80%"
$formatDonor.Copy() | Out-Null
$ws.Range("H18").PasteSpecial(-4122) | Out-Null
$ws.Range("K18").Value = "10.0 MJ/m2"
$ws.Range("O18").Value = "10.0 °C"
$ws.Range("E19").Value = "2026-02-17 16:19:32"
$ws.Range("K19").Value = "9.2 MJ/m2"
$ws.Range("O19").Value = "7.2 °C"
$ws.Range("E20").Value = "2026-02-17 16:19:35"
$ws.Range("H20").Value = "'56%"
$formatDonor.Copy() | Out-Null
$ws.Range("H20").PasteSpecial(-4122) | Out-Null
$ws.Range("K20").Value = "9.4 MJ/m2"
$ws.Range("O20").Value = "-2.2 °C"
$ws.Range("E21").Value = "2026-02-17 16:19:37"
$ws.Range("K21").Value = "6.0 MJ/m2"
$ws.Range("O21").Value = "9.6 °C"
$ws.Range("E22").Value = "2026-02-17 16:19:40"
$ws.Range("E23").Value = "2026-02-17 16:19:43"
$ws.Range("I23").Value = "2.3 mm"
$ws.Range("K23").Value = "11.6 MJ/m2"
$ws.Range("O23").Value = "-4.2 °C"
$ws.Range("E24").Value = "2026-02-17 16:19:45"
$ws.Range("K24").Value = "14.4 MJ/m2"
$ws.Range("O24").Value = "12.4 °C"
$ws.Range("E25").Value = "2026-02-17 16:19:48"
$ws.Range("H25").Value = "'48%"
$formatDonor.Copy() | Out-Null
$ws.Range("H25").PasteSpecial(-4122) | Out-Null
$ws.Range("K25").Value = "13.2 MJ/m2"
$ws.Range("O25").Value = "-1.2 °C"
$ws.Range("E26").Value = "2026-02-17 16:19:51"
$ws.Range("E27").Value = "2026-02-17 16:19:53"
$ws.Range("K27").Value = "10.1 MJ/m2"
$ws.Range("O27").Value = "-0.7 °C"
$ws.Range("E28").Value = "2026-02-17 16:19:56"
$ws.Range("K28").Value = "9.3 MJ/m2"
$ws.Range("O28").Value = "8.5 °C"
$ws.Range("E29").Value = "2026-02-17 16:19:58"
$ws.Range("E30").Value = "2026-02-17 16:20:01"
$ws.Range("J30").Value = "1017.9 hPa"
$ws.Range("K30").Value = "7.7 MJ/m2"
$ws.Range("E31").Value = "2026-02-17 16:20:03"
$ws.Range("K31").Value = "9.2 MJ/m2"
$ws.Range("E32").Value = "2026-02-17 16:20:06"
$ws.Range("H32").Value = "'69%"
$formatDonor.Copy() | Out-Null
$ws.Range("H32").PasteSpecial(-4122) | Out-Null
$ws.Range("K32").Value = "10.5 MJ/m2"
$ws.Range("O32").Value = "8.2 °C"
$ws.Range("E33").Value = "2026-02-17 16:20:09"
$ws.Range("H33").Value = "'38%"
$formatDonor.Copy() | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null
$ws.Range("K33").Value = "6.8 MJ/m2"
$ws.Range("O33").Value = "6.4 °C"
$ws.Range("E34").Value = "2026-02-17 16:20:11"
$ws.Range("K34").Value = "10.9 MJ/m2"
$ws.Range("E35").Value = "2026-02-17 16:20:14"
$ws.Range("K35").Value = "8.8 MJ/m2"
$ws.Range("O35").Value = "7.0 °C"
$ws.Range("E36").Value = "2026-02-17 16:20:17"
$ws.Range("K36").Value = "10.4 MJ/m2"
$ws.Range("E37").Value = "2026-02-17 16:20:19"
$ws.Range("J37").Value = "1018.5 hPa"
$ws.Range("E38").Value = "2026-02-17 16:20:22"
$ws.Range("K38").Value = "10.4 MJ/m2"
$ws.Range("O38").Value = "11.3 °C"
$ws.Range("E39").Value = "2026-02-17 16:20:24"
$ws.Range("E40").Value = "2026-02-17 16:20:27"
$ws.Range("O40").Value = "9.5 °C"
$ws.Range("E41").Value = "2026-02-17 16:20:29"
$ws.Range("H41").Value = "'48%"
$formatDonor.Copy() | Out-Null
$ws.Range("H41").PasteSpecial(-4122) | Out-Null
$ws.Range("J41").Value = "1017.9 hPa"
$ws.Range("K41").Value = "11.8 MJ/m2"
$ws.Range("M41").Value = "22.1 °C 15:35 TU"
$ws.Range("O41").Value = "16.5 °C"
$ws.Range("E42").Value = "2026-02-17 16:20:32"
$ws.Range("O42").Value = "13.2 °C"
$ws.Range("E43").Value = "2026-02-17 16:20:35"
$ws.Range("K43").Value = "12.9 MJ/m2"
$ws.Range("L43").Value = "29.2 km/h - 188º 15:35 TU"
$ws.Range("O43").Value = "7.7 °C"
$ws.Range("E44").Value = "2026-02-17 16:20:37"
$ws.Range("K44").Value = "9.6 MJ/m2"
$ws.Range("O44").Value = "-3.5 °C"
$ws.Range("E45").Value = "2026-02-17 16:20:40"
$ws.Range("H45").Value = "'63%"
$formatDonor.Copy() | Out-Null
$ws.Range("H45").PasteSpecial(-4122) | Out-Null
$ws.Range("J45").Value = "1021.5 hPa"
$ws.Range("K45").Value = "3.3 MJ/m2"
$ws.Range("O45").Value = "5.5 °C"
$ws.Range("E46").Value = "2026-02-17 16:20:42"
$ws.Range("K46").Value = "13.7 MJ/m2"
$ws.Range("O46").Value = "15.5 °C"

$excel.CutCopyMode = 0
